# The commit renames three embedded pictures' internal OOXML "name"
# attributes (the wp:docPr / pic:cNvPr pair for each <w:drawing>):
#   - footer1.xml picture (id=1, Pearson logo): image1.png -> image2.png
#   - footer2.xml picture (id=2, Pearson logo): image1.png -> image2.png
#   - header2.xml picture (id=3, BTec logo):    image2.jpg -> image1.jpg
#
# InlineShape objects in the Word object model do not expose a settable
# "Name" property (that only exists on floating Shape objects), so the
# rename is applied through the document's raw OOXML package, which is a
# supported, documented Word COM surface (Document.WordOpenXML).

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

# Pearson Edexcel logo pictures (two occurrences, id="1" and id="2"),
# both currently named image1.png -> rename to image2.png.
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="1" name="image1.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="1" name="image2.png"')
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"')
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"')

# BTec logo picture (id="3"), currently named image2.jpg -> rename to image1.jpg.
$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="3" name="image2.jpg"',
    'descr="BTec_Logo-Orange" id="3" name="image1.jpg"')
$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="0" name="image2.jpg"',
    'descr="BTec_Logo-Orange" id="0" name="image1.jpg"')

$d.WordOpenXML = $xml

Write-Host "Renamed inline picture names in the document OOXML."
